# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Wed Nov 15 10:48:08 UTC 2023 with GitHub Actions"
#
# Price (column D) cells are plain-text strings formatted like "1.234.56" or
# "245.30" that LOOK numeric. A bare Range.Value assignment would let Excel
# auto-convert those into real numbers (dropping the trailing zero, turning
# "1.234.56" into a date, flipping the cell style to a text-format xf, etc).
# Set-CryptoText forces the cell to Text format just long enough to write the
# literal string, then restores the cell's original Style so no formatting
# diff leaks in beyond the value itself.
function Set-CryptoText($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-CryptoText $ws.Range("D2") '35.836.16'
$ws.Range("E2").Value = '  -1.70%  '

# Row 3
Set-CryptoText $ws.Range("D3") '1.986.72'
$ws.Range("E3").Value = '  -2.74%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
Set-CryptoText $ws.Range("D5") '245.30'
$ws.Range("E5").Value = '  -0.29%  '

# Row 6
$ws.Range("E6").Value = '  -3.00%  '

# Row 7
Set-CryptoText $ws.Range("D7") '59.96'
$ws.Range("E7").Value = '  +11.33%  '

# Row 8
$ws.Range("E8").Value = '  -0.09%  '

# Row 9
Set-CryptoText $ws.Range("D9") '58.10'
$ws.Range("E9").Value = '  -6.48%  '

# Row 10
$ws.Range("E10").Value = '  +0.55%  '

# Row 11
Set-CryptoText $ws.Range("D11") '0.0740'
$ws.Range("E11").Value = '  -0.60%  '

# Row 12
$ws.Range("E12").Value = '  -2.15%  '

# Row 13
$ws.Range("E13").Value = '  +0.36%  '

# Row 14
Set-CryptoText $ws.Range("D14") '14.74'
$ws.Range("E14").Value = '  +1.36%  '

# Row 15
$ws.Range("E15").Value = '  -3.00%  '

# Row 16
$ws.Range("E16").Value = '  -0.60%  '

# Row 17
Set-CryptoText $ws.Range("D17") '18.98'
$ws.Range("E17").Value = '  +12.42%  '

# Row 18
Set-CryptoText $ws.Range("D18") '1.985.27'
$ws.Range("E18").Value = '  -2.92%  '

# Row 19
Set-CryptoText $ws.Range("D19") '35.773.53'
$ws.Range("E19").Value = '  -1.58%  '

# Row 20
Set-CryptoText $ws.Range("D20") '71.56'
$ws.Range("E20").Value = '  +0.19%  '

# Row 21
Set-CryptoText $ws.Range("D21") '0.0₃0849'
$ws.Range("E21").Value = '  -0.38%  '

# Row 22
Set-CryptoText $ws.Range("D22") '5.20'
$ws.Range("E22").Value = '  +0.85%  '

# Row 23
Set-CryptoText $ws.Range("D23") '232.75'
$ws.Range("E23").Value = '  -1.54%  '

# Row 24
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
Set-CryptoText $ws.Range("D25") '2.58'
$ws.Range("E25").Value = '  +15.02%  '

# Row 26
Set-CryptoText $ws.Range("D26") '2.28'
$ws.Range("E26").Value = '  -4.34%  '

# Row 27
Set-CryptoText $ws.Range("D27") '9.62'
$ws.Range("E27").Value = '  +5.65%  '

# Row 28
Set-CryptoText $ws.Range("D28") '165.30'
$ws.Range("E28").Value = '  +0.50%  '

# Row 29
Set-CryptoText $ws.Range("D29") '19.37'
$ws.Range("E29").Value = '  -2.05%  '

# Row 30
$ws.Range("E30").Value = '  -0.59%  '

# Row 31
$ws.Range("E31").Value = '  +0.36%  '

# Row 32
$ws.Range("E32").Value = '  -2.89%  '

# Row 33
$ws.Range("E33").Value = '  +13.19%  '

# Row 34
Set-CryptoText $ws.Range("D34") '0.0600'
$ws.Range("E34").Value = '  +1.99%  '

# Row 35
$ws.Range("E35").Value = '  +0.87%  '

# Row 36
$ws.Range("E36").Value = '  +11.26%  '

# Row 38
$ws.Range("E38").Value = '  -2.25%  '

# Row 39
Set-CryptoText $ws.Range("D39") '5.69'
$ws.Range("E39").Value = '  +13.72%  '

# Row 40
$ws.Range("E40").Value = '  +0.19%  '

# Row 42
$ws.Range("E42").Value = '  +5.10%  '

# Row 43
$ws.Range("E43").Value = '  +0.87%  '

# Row 44
Set-CryptoText $ws.Range("D44") '1.11'
$ws.Range("E44").Value = '  +1.02%  '

# Row 45
Set-CryptoText $ws.Range("D45") '16.56'
$ws.Range("E45").Value = '  +5.90%  '

# Row 46
Set-CryptoText $ws.Range("D46") '7.77'
$ws.Range("E46").Value = '  +5.62%  '

# Row 47
Set-CryptoText $ws.Range("D47") '93.15'
$ws.Range("E47").Value = '  -0.22%  '

# Row 48
Set-CryptoText $ws.Range("D48") '1.360.76'
$ws.Range("E48").Value = '  -2.62%  '

# Row 49
Set-CryptoText $ws.Range("D49") '2.90'
$ws.Range("E49").Value = '  -0.32%  '

# Row 50
$ws.Range("E50").Value = '  +2.73%  '

# Row 51
Set-CryptoText $ws.Range("D51") '46.58'
$ws.Range("E51").Value = '  +3.58%  '
